# Commit: Fri, May 29, 2020  9:05:09 AM
#
# The table on slide 16 (the "PLENARY- COMPLETE THE MISSING GAPS" slide)
# has its table style switched from the deck's default table style
# ({CA177CB4-AF6C-461D-B295-C07ECAF9CF51}) to a different built-in table
# style ({B7D86948-5166-413B-9F71-2671AAC12FC5}), leaving the banding /
# first-row / first-column flags untouched.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(16)
$sh = $s.Shapes.Item(3)
$tbl = $sh.Table

$tbl.ApplyStyle("{B7D86948-5166-413B-9F71-2671AAC12FC5}", $true)
